$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: clear the cluster/gene name columns (A:D) for existing rows so their
# shared-string entries are dropped from the table, then rebuild the table in
# the exact order required by the target file (ECs, FAPs, Inflammatory-Mac,
# MuSCs, Resolving-Mac, Btla, Cd79a) before writing the refreshed data grid.
for ($r = 2; $r -le 21; $r++) {
    $ws.Range("A" + $r).Value = ""
    $ws.Range("B" + $r).Value = ""
    $ws.Range("C" + $r).Value = ""
    $ws.Range("D" + $r).Value = ""
}

# Step 2: seed a scratch cell with each cluster/gene name in the desired order
# so the shared-string table is rebuilt with that exact ordering.
$seedCell = $ws.Range("V1")
$seedCell.Value = "ECs"
$seedCell.Value = "FAPs"
$seedCell.Value = "Inflammatory-Mac"
$seedCell.Value = "MuSCs"
$seedCell.Value = "Resolving-Mac"
$seedCell.Value = "Btla"
$seedCell.Value = "Cd79a"
$seedCell.Value = ""

# Step 3: write the refreshed data grid (rows 2-26).
# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Btla"
$ws.Range("C2").Value = "Cd79a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.04270366666666667
$ws.Range("H2").Value = 0.128111
$ws.Range("I2").Value = 0.006976065356144797
$ws.Range("J2").Value = 0.007176859668027315
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.2905836666666666
$ws.Range("N2").Value = 0.8717509999999999
$ws.Range("O2").Value = 0.0857173605651528
$ws.Range("P2").Value = 0.09660500428308474
$ws.Range("Q2").Value = 0.01240898804011111
$ws.Range("R2").Value = 0.111680892361
$ws.Range("S2").Value = 0.0005979699094587347
$ws.Range("T2").Value = 0.0006933205589688769

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Btla"
$ws.Range("C3").Value = "Cd79a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.04270366666666667
$ws.Range("H3").Value = 0.128111
$ws.Range("I3").Value = 0.006976065356144797
$ws.Range("J3").Value = 0.007176859668027315
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.535638666666667
$ws.Range("N3").Value = 4.606916
$ws.Range("O3").Value = 0.4529879287381047
$ws.Range("P3").Value = 0.5105255284041104
$ws.Range("Q3").Value = 0.06557740174177777
$ws.Range("R3").Value = 0.590196615676
$ws.Range("S3").Value = 0.00316007339642168
$ws.Range("T3").Value = 0.003663970074301793

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Btla"
$ws.Range("C4").Value = "Cd79a"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.04270366666666667
$ws.Range("H4").Value = 0.128111
$ws.Range("I4").Value = 0.006976065356144797
$ws.Range("J4").Value = 0.007176859668027315
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3897316666666666
$ws.Range("N4").Value = 1.169195
$ws.Range("O4").Value = 0.1149643755911652
$ws.Range("P4").Value = 0.1295669153035228
$ws.Range("Q4").Value = 0.01664297118277778
$ws.Range("R4").Value = 0.149786740645
$ws.Range("S4").Value = 0.0008019989977523459
$ws.Range("T4").Value = 0.0009298835687525637

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Btla"
$ws.Range("C5").Value = "Cd79a"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.04270366666666667
$ws.Range("H5").Value = 0.128111
$ws.Range("I5").Value = 0.006976065356144797
$ws.Range("J5").Value = 0.007176859668027315
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.1461935
$ws.Range("N5").Value = 2.292387
$ws.Range("O5").Value = 0.3381080658935905
$ws.Range("P5").Value = 0.2540359069889083
$ws.Range("Q5").Value = 0.0489466651595
$ws.Range("R5").Value = 0.293679990957
$ws.Range("S5").Value = 0.002358663965113399
$ws.Range("T5").Value = 0.001823180055099434

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Btla"
$ws.Range("C6").Value = "Cd79a"
$ws.Range("D6").Value = "Resolving-Mac"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.04270366666666667
$ws.Range("H6").Value = 0.128111
$ws.Range("I6").Value = 0.006976065356144797
$ws.Range("J6").Value = 0.007176859668027315
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02787366666666667
$ws.Range("N6").Value = 0.083621
$ws.Range("O6").Value = 0.00822226921198673
$ws.Range("P6").Value = 0.009266645020373741
$ws.Range("Q6").Value = 0.001190307770111111
$ws.Range("R6").Value = 0.010712769931
$ws.Range("S6").Value = 0.00005735908739864
$ws.Range("T6").Value = 0.00006650541090465

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Btla"
$ws.Range("C7").Value = "Cd79a"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.07261366666666667
$ws.Range("H7").Value = 0.217841
$ws.Range("I7").Value = 0.01186215901248089
$ws.Range("J7").Value = 0.01220359131489676
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.2905836666666666
$ws.Range("N7").Value = 0.8717509999999999
$ws.Range("O7").Value = 0.0857173605651528
$ws.Range("P7").Value = 0.09660500428308474
$ws.Range("Q7").Value = 0.02110034551011111
$ws.Range("R7").Value = 0.189903109591
$ws.Range("S7").Value = 0.001016792961154001
$ws.Range("T7").Value = 0.001178927991244617

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Btla"
$ws.Range("C8").Value = "Cd79a"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.07261366666666667
$ws.Range("H8").Value = 0.217841
$ws.Range("I8").Value = 0.01186215901248089
$ws.Range("J8").Value = 0.01220359131489676
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.535638666666667
$ws.Range("N8").Value = 4.606916
$ws.Range("O8").Value = 0.4529879287381047
$ws.Range("P8").Value = 0.5105255284041104
$ws.Range("Q8").Value = 0.1115083542617778
$ws.Range("R8").Value = 1.003575188356
$ws.Range("S8").Value = 0.005373414841425759
$ws.Range("T8").Value = 0.006230244904465478

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Btla"
$ws.Range("C9").Value = "Cd79a"
$ws.Range("D9").Value = "Inflammatory-Mac"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.07261366666666667
$ws.Range("H9").Value = 0.217841
$ws.Range("I9").Value = 0.01186215901248089
$ws.Range("J9").Value = 0.01220359131489676
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.3897316666666666
$ws.Range("N9").Value = 1.169195
$ws.Range("O9").Value = 0.1149643755911652
$ws.Range("P9").Value = 0.1295669153035228
$ws.Range("Q9").Value = 0.02829984533277778
$ws.Range("R9").Value = 0.254698607995
$ws.Range("S9").Value = 0.001363725704032978
$ws.Range("T9").Value = 0.001581181682296034

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Btla"
$ws.Range("C10").Value = "Cd79a"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.07261366666666667
$ws.Range("H10").Value = 0.217841
$ws.Range("I10").Value = 0.01186215901248089
$ws.Range("J10").Value = 0.01220359131489676
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.1461935
$ws.Range("N10").Value = 2.292387
$ws.Range("O10").Value = 0.3381080658935905
$ws.Range("P10").Value = 0.2540359069889083
$ws.Range("Q10").Value = 0.0832293127445
$ws.Range("R10").Value = 0.4993758764669999
$ws.Range("S10").Value = 0.004010691641032136
$ws.Range("T10").Value = 0.003100150388201761

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Btla"
$ws.Range("C11").Value = "Cd79a"
$ws.Range("D11").Value = "Resolving-Mac"
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.07261366666666667
$ws.Range("H11").Value = 0.217841
$ws.Range("I11").Value = 0.01186215901248089
$ws.Range("J11").Value = 0.01220359131489676
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.02787366666666667
$ws.Range("N11").Value = 0.083621
$ws.Range("O11").Value = 0.00822226921198673
$ws.Range("P11").Value = 0.009266645020373741
$ws.Range("Q11").Value = 0.002024009140111111
$ws.Range("R11").Value = 0.018216082261
$ws.Range("S11").Value = 0.00009753386483601
$ws.Range("T11").Value = 0.0001130863486888643

# Row 12
$ws.Range("A12").Value = "Inflammatory-Mac"
$ws.Range("B12").Value = "Btla"
$ws.Range("C12").Value = "Cd79a"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.695177999999999
$ws.Range("H12").Value = 14.085534
$ws.Range("I12").Value = 0.7670036590160066
$ws.Range("J12").Value = 0.7890805697186616
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.2905836666666666
$ws.Range("N12").Value = 0.8717509999999999
$ws.Range("O12").Value = 0.0857173605651528
$ws.Range("P12").Value = 0.09660500428308474
$ws.Range("Q12").Value = 1.364342038892666
$ws.Range("R12").Value = 12.279078350034
$ws.Range("S12").Value = 0.06574552919466656
$ws.Range("T12").Value = 0.07622913181737025

# Row 13
$ws.Range("A13").Value = "Inflammatory-Mac"
$ws.Range("B13").Value = "Btla"
$ws.Range("C13").Value = "Cd79a"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.695177999999999
$ws.Range("H13").Value = 14.085534
$ws.Range("I13").Value = 0.7670036590160066
$ws.Range("J13").Value = 0.7890805697186616
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.535638666666667
$ws.Range("N13").Value = 4.606916
$ws.Range("O13").Value = 0.4529879287381047
$ws.Range("P13").Value = 0.5105255284041104
$ws.Range("Q13").Value = 7.210096883682666
$ws.Range("R13").Value = 64.890871953144
$ws.Range("S13").Value = 0.3474433988322084
$ws.Range("T13").Value = 0.4028457748090362

# Row 14
$ws.Range("A14").Value = "Inflammatory-Mac"
$ws.Range("B14").Value = "Btla"
$ws.Range("C14").Value = "Cd79a"
$ws.Range("D14").Value = "Inflammatory-Mac"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 4.695177999999999
$ws.Range("H14").Value = 14.085534
$ws.Range("I14").Value = 0.7670036590160066
$ws.Range("J14").Value = 0.7890805697186616
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 0.3897316666666666
$ws.Range("N14").Value = 1.169195
$ws.Range("O14").Value = 0.1149643755911652
$ws.Range("P14").Value = 0.1295669153035228
$ws.Range("Q14").Value = 1.829859547236666
$ws.Range("R14").Value = 16.46873592513
$ws.Range("S14").Value = 0.08817809673491417
$ws.Range("T14").Value = 0.1022387353443933

# Row 15
$ws.Range("A15").Value = "Inflammatory-Mac"
$ws.Range("B15").Value = "Btla"
$ws.Range("C15").Value = "Cd79a"
$ws.Range("D15").Value = "MuSCs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 4.695177999999999
$ws.Range("H15").Value = 14.085534
$ws.Range("I15").Value = 0.7670036590160066
$ws.Range("J15").Value = 0.7890805697186616
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 1.1461935
$ws.Range("N15").Value = 2.292387
$ws.Range("O15").Value = 0.3381080658935905
$ws.Range("P15").Value = 0.2540359069889083
$ws.Range("Q15").Value = 5.381582504942998
$ws.Range("R15").Value = 32.28949502965799
$ws.Range("S15").Value = 0.259330123683209
$ws.Range("T15").Value = 0.2004547982158047

# Row 16
$ws.Range("A16").Value = "Inflammatory-Mac"
$ws.Range("B16").Value = "Btla"
$ws.Range("C16").Value = "Cd79a"
$ws.Range("D16").Value = "Resolving-Mac"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 4.695177999999999
$ws.Range("H16").Value = 14.085534
$ws.Range("I16").Value = 0.7670036590160066
$ws.Range("J16").Value = 0.7890805697186616
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.02787366666666667
$ws.Range("N16").Value = 0.083621
$ws.Range("O16").Value = 0.00822226921198673
$ws.Range("P16").Value = 0.009266645020373741
$ws.Range("Q16").Value = 0.1308718265126667
$ws.Range("R16").Value = 1.177846438614
$ws.Range("S16").Value = 0.006306510571008479
$ws.Range("T16").Value = 0.00731212953205711

# Row 17
$ws.Range("A17").Value = "MuSCs"
$ws.Range("B17").Value = "Btla"
$ws.Range("C17").Value = "Cd79a"
$ws.Range("D17").Value = "ECs"
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.5137985
$ws.Range("H17").Value = 1.027597
$ws.Range("I17").Value = 0.08393405521514537
$ws.Range("J17").Value = 0.05756663724649613
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.2905836666666666
$ws.Range("N17").Value = 0.8717509999999999
$ws.Range("O17").Value = 0.0857173605651528
$ws.Range("P17").Value = 0.09660500428308474
$ws.Range("Q17").Value = 0.1493014520578333
$ws.Range("R17").Value = 0.895808712347
$ws.Range("S17").Value = 0.007194605674572059
$ws.Range("T17").Value = 0.005561225237760544

# Row 18
$ws.Range("A18").Value = "MuSCs"
$ws.Range("B18").Value = "Btla"
$ws.Range("C18").Value = "Cd79a"
$ws.Range("D18").Value = "FAPs"
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.5137985
$ws.Range("H18").Value = 1.027597
$ws.Range("I18").Value = 0.08393405521514537
$ws.Range("J18").Value = 0.05756663724649613
$ws.Range("K18").Value = 3
$ws.Range("L18").Value = 1
$ws.Range("M18").Value = 1.535638666666667
$ws.Range("N18").Value = 4.606916
$ws.Range("O18").Value = 0.4529879287381047
$ws.Range("P18").Value = 0.5105255284041104
$ws.Range("Q18").Value = 0.7890088434753334
$ws.Range("R18").Value = 4.734053060852
$ws.Range("S18").Value = 0.03802111382249842
$ws.Range("T18").Value = 0.02938923789871517

# Row 19
$ws.Range("A19").Value = "MuSCs"
$ws.Range("B19").Value = "Btla"
$ws.Range("C19").Value = "Cd79a"
$ws.Range("D19").Value = "Inflammatory-Mac"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.5137985
$ws.Range("H19").Value = 1.027597
$ws.Range("I19").Value = 0.08393405521514537
$ws.Range("J19").Value = 0.05756663724649613
$ws.Range("K19").Value = 3
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 0.3897316666666666
$ws.Range("N19").Value = 1.169195
$ws.Range("O19").Value = 0.1149643755911652
$ws.Range("P19").Value = 0.1295669153035228
$ws.Range("Q19").Value = 0.2002435457358333
$ws.Range("R19").Value = 1.201461274415
$ws.Range("S19").Value = 0.009649426248643568
$ws.Range("T19").Value = 0.007458731612425383

# Row 20
$ws.Range("A20").Value = "MuSCs"
$ws.Range("B20").Value = "Btla"
$ws.Range("C20").Value = "Cd79a"
$ws.Range("D20").Value = "MuSCs"
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.5137985
$ws.Range("H20").Value = 1.027597
$ws.Range("I20").Value = 0.08393405521514537
$ws.Range("J20").Value = 0.05756663724649613
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 1.1461935
$ws.Range("N20").Value = 2.292387
$ws.Range("O20").Value = 0.3381080658935905
$ws.Range("P20").Value = 0.2540359069889083
$ws.Range("Q20").Value = 0.58891250100975
$ws.Range("R20").Value = 2.355650004039
$ws.Range("S20").Value = 0.02837878107139863
$ws.Range("T20").Value = 0.01462399290521511

# Row 21
$ws.Range("A21").Value = "MuSCs"
$ws.Range("B21").Value = "Btla"
$ws.Range("C21").Value = "Cd79a"
$ws.Range("D21").Value = "Resolving-Mac"
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.5137985
$ws.Range("H21").Value = 1.027597
$ws.Range("I21").Value = 0.08393405521514537
$ws.Range("J21").Value = 0.05756663724649613
$ws.Range("K21").Value = 1
$ws.Range("L21").Value = 0.3333333333333333
$ws.Range("M21").Value = 0.02787366666666667
$ws.Range("N21").Value = 0.083621
$ws.Range("O21").Value = 0.00822226921198673
$ws.Range("P21").Value = 0.009266645020373741
$ws.Range("Q21").Value = 0.01432144812283334
$ws.Range("R21").Value = 0.085928688737
$ws.Range("S21").Value = 0.0006901283980326839
$ws.Range("T21").Value = 0.0005334495923799049

# Row 22
$ws.Range("A22").Value = "Resolving-Mac"
$ws.Range("B22").Value = "Btla"
$ws.Range("C22").Value = "Cd79a"
$ws.Range("D22").Value = "ECs"
$ws.Range("E22").Value = 2
$ws.Range("F22").Value = 0.6666666666666666
$ws.Range("G22").Value = 0.7971606666666666
$ws.Range("H22").Value = 2.391482
$ws.Range("I22").Value = 0.1302240614002222
$ws.Range("J22").Value = 0.1339723420519182
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.2905836666666666
$ws.Range("N22").Value = 0.8717509999999999
$ws.Range("O22").Value = 0.0857173605651528
$ws.Range("P22").Value = 0.09660500428308474
$ws.Range("Q22").Value = 0.2316418694424444
$ws.Range("R22").Value = 2.084776824982
$ws.Range("S22").Value = 0.01116246282530144
$ws.Range("T22").Value = 0.01294239867774046

# Row 23
$ws.Range("A23").Value = "Resolving-Mac"
$ws.Range("B23").Value = "Btla"
$ws.Range("C23").Value = "Cd79a"
$ws.Range("D23").Value = "FAPs"
$ws.Range("E23").Value = 2
$ws.Range("F23").Value = 0.6666666666666666
$ws.Range("G23").Value = 0.7971606666666666
$ws.Range("H23").Value = 2.391482
$ws.Range("I23").Value = 0.1302240614002222
$ws.Range("J23").Value = 0.1339723420519182
$ws.Range("K23").Value = 3
$ws.Range("L23").Value = 1
$ws.Range("M23").Value = 1.535638666666667
$ws.Range("N23").Value = 4.606916
$ws.Range("O23").Value = 0.4529879287381047
$ws.Range("P23").Value = 0.5105255284041104
$ws.Range("Q23").Value = 1.224150743279111
$ws.Range("R23").Value = 11.017356689512
$ws.Range("S23").Value = 0.05898992784555044
$ws.Range("T23").Value = 0.06839630071759177

# Row 24
$ws.Range("A24").Value = "Resolving-Mac"
$ws.Range("B24").Value = "Btla"
$ws.Range("C24").Value = "Cd79a"
$ws.Range("D24").Value = "Inflammatory-Mac"
$ws.Range("E24").Value = 2
$ws.Range("F24").Value = 0.6666666666666666
$ws.Range("G24").Value = 0.7971606666666666
$ws.Range("H24").Value = 2.391482
$ws.Range("I24").Value = 0.1302240614002222
$ws.Range("J24").Value = 0.1339723420519182
$ws.Range("K24").Value = 3
$ws.Range("L24").Value = 1
$ws.Range("M24").Value = 0.3897316666666666
$ws.Range("N24").Value = 1.169195
$ws.Range("O24").Value = 0.1149643755911652
$ws.Range("P24").Value = 0.1295669153035228
$ws.Range("Q24").Value = 0.3106787552211111
$ws.Range("R24").Value = 2.79610879699
$ws.Range("S24").Value = 0.0149711279058221
$ws.Range("T24").Value = 0.01735838309565547

# Row 25
$ws.Range("A25").Value = "Resolving-Mac"
$ws.Range("B25").Value = "Btla"
$ws.Range("C25").Value = "Cd79a"
$ws.Range("D25").Value = "MuSCs"
$ws.Range("E25").Value = 2
$ws.Range("F25").Value = 0.6666666666666666
$ws.Range("G25").Value = 0.7971606666666666
$ws.Range("H25").Value = 2.391482
$ws.Range("I25").Value = 0.1302240614002222
$ws.Range("J25").Value = 0.1339723420519182
$ws.Range("K25").Value = 2
$ws.Range("L25").Value = 1
$ws.Range("M25").Value = 1.1461935
$ws.Range("N25").Value = 2.292387
$ws.Range("O25").Value = 0.3381080658935905
$ws.Range("P25").Value = 0.2540359069889083
$ws.Range("Q25").Value = 0.9137003745889999
$ws.Range("R25").Value = 5.482202247533999
$ws.Range("S25").Value = 0.04402980553283731
$ws.Range("T25").Value = 0.03403378542458731

# Row 26
$ws.Range("A26").Value = "Resolving-Mac"
$ws.Range("B26").Value = "Btla"
$ws.Range("C26").Value = "Cd79a"
$ws.Range("D26").Value = "Resolving-Mac"
$ws.Range("E26").Value = 2
$ws.Range("F26").Value = 0.6666666666666666
$ws.Range("G26").Value = 0.7971606666666666
$ws.Range("H26").Value = 2.391482
$ws.Range("I26").Value = 0.1302240614002222
$ws.Range("J26").Value = 0.1339723420519182
$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 0.3333333333333333
$ws.Range("M26").Value = 0.02787366666666667
$ws.Range("N26").Value = 0.083621
$ws.Range("O26").Value = 0.00822226921198673
$ws.Range("P26").Value = 0.009266645020373741
$ws.Range("Q26").Value = 0.02221979070244445
$ws.Range("R26").Value = 0.199978116322
$ws.Range("S26").Value = 0.001070737290710917
$ws.Range("T26").Value = 0.001241474136343216

Write-Host "edit applied"
